$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 171, shifting existing rows 171:197 down to 172:198.
# The new row inherits formatting (number format, etc.) from the row above.
$ws.Rows("171:171").Insert()

# Populate the new row 171 with the latest weekly entry.
$ws.Cells.Item(171, 1).Value = 7
$ws.Cells.Item(171, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(171, 3).Value = "Ñuble"
$ws.Cells.Item(171, 4).Value = 44505
$ws.Cells.Item(171, 5).Value = 16
$ws.Cells.Item(171, 6).Value = 100114013
$ws.Cells.Item(171, 7).Value = "Zanahoria"
$ws.Cells.Item(171, 8).Value = "Sin especificar"
$ws.Cells.Item(171, 9).Value = "Primera"
$ws.Cells.Item(171, 10).Value = 120
$ws.Cells.Item(171, 11).Value = 7500
$ws.Cells.Item(171, 12).Value = 8000
$ws.Cells.Item(171, 13).Value = 7750
$ws.Cells.Item(171, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(171, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(171, 16).Value = 388
$ws.Cells.Item(171, 17).Value = 20
$ws.Cells.Item(171, 18).Value = "Hortaliza"
